# checkpoint grupo 13; feito g13.5b
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:D1) formatting: add thin border all around and
#     set vertical alignment to top (keeps existing bold font + centered text) ---
$headerRng = $ws.Range("A1:D1")
$headerRng.Borders.LineStyle = 1
$headerRng.VerticalAlignment = -4160

# --- Update the "Variável" column label for every data row ---
$ws.Range("B2:B10").Value = "Diferença 2024/07 - 2024/07"

# --- Update the "Valor" column values ---
$ws.Range("C2").Value = 1.19
$ws.Range("C3").Value = 1.13
$ws.Range("C4").Value = 0.99
$ws.Range("C5").Value = 0.96
$ws.Range("C6").Value = 0.9399999999999999
$ws.Range("C7").Value = 0.93
$ws.Range("C8").Value = 0.86
$ws.Range("C9").Value = 0.83
$ws.Range("C10").Value = 0.8100000000000001

# --- Rows reshuffled: region names change (ranking refresh) ---
$ws.Range("A6").Value = "Rio Grande do Sul"
$ws.Range("A7").Value = "Tocantins"
$ws.Range("A9").Value = "Brasil"
$ws.Range("A10").Value = "Nordeste"

# --- Updated placement ("Colocação") for Sergipe ---
$ws.Range("D8").Value = "12º"

# --- Page margins reset to Excel's standard "Normal" preset ---
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72
